$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Row 33: note row ("Improved state machine -> ...")
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = "Improved state machine -> Waits for another handshake after end of transmission and exits handshake if too much time between transitions"

# ---------------------------------------------------------------------------
# Row 34: new data row (own, non-shared H/I formulas; G34 gets a formula of
# its own, kept apart from the G35:G40 fill-down group below)
# ---------------------------------------------------------------------------
$ws.Range("A34").Value = 15
$ws.Range("B34").Value = 150
$ws.Range("C34").Value = 10
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("J34").Value = 10

$ws.Range("G34:G40").Formula = '=IF(F34,E34/F34,"")'
$ws.Range("H34").Formula = "=118-F34"
$ws.Range("I34").Formula = "=H34/118"

# ---------------------------------------------------------------------------
# Rows 35-40: fill-down block, H/I/G sharing formulas across the range.
# Row 37 is the "Reset more variables..." note row and keeps the G formula
# but loses its H/I formulas (matching the pattern already used by row 18).
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = 15
$ws.Range("B35").Value = 150
$ws.Range("C35").Value = 10
$ws.Range("E35").Value = 6
$ws.Range("F35").Value = 12
$ws.Range("J35").Value = 10

$ws.Range("A36").Value = 15
$ws.Range("B36").Value = 150
$ws.Range("C36").Value = 10
$ws.Range("E36").Value = 4
$ws.Range("F36").Value = 12
$ws.Range("J36").Value = 10

$ws.Range("A37").Value = "Reset more variables between transmissions, most importantly ""bits"""

$ws.Range("A38").Value = 15
$ws.Range("B38").Value = 150
$ws.Range("C38").Value = 10
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("J38").Value = 10

$ws.Range("A39").Value = 15
$ws.Range("B39").Value = 150
$ws.Range("C39").Value = 10
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("J39").Value = 10

$ws.Range("A40").Value = 15
$ws.Range("B40").Value = 150
$ws.Range("C40").Value = 10
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("J40").Value = 10

# Fill-down shared formulas for G35:G40, H35:H40, I35:I40
$ws.Range("G35:G40").Formula = '=IF(F35,E35/F35,"")'
$ws.Range("H35:H40").Formula = "=118-F35"
$ws.Range("I35:I40").Formula = "=H35/118"

# Row 37 is a note row: it keeps the G formula, but H/I are cleared back to
# blank (still formatted) cells, same treatment as row 18's H/I.
$ws.Range("H37:I37").ClearContents()

# ---------------------------------------------------------------------------
# Row 41: note row ("Handshake is now based on...")
# ---------------------------------------------------------------------------
$ws.Range("A41").Value = "Handshake is now based on initial state of the bulb"

# ---------------------------------------------------------------------------
# Number formats for the newly-written G/H/I cells (reuses the workbook's
# existing "Percent"/"0" styles rather than creating new ones).
# ---------------------------------------------------------------------------
$ws.Range("G34:G40").NumberFormat = "0.0%"
$ws.Range("H34:H40").NumberFormat = "0"
$ws.Range("I34:I40").NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Also clear the old H18/I18 formulas (transmission considered over once the
# bulb stays on => that earlier sample no longer reports a handshake count).
# ---------------------------------------------------------------------------
$ws.Range("H18:I18").ClearContents()

# ---------------------------------------------------------------------------
# Final selection, matching the state the workbook was left in.
# ---------------------------------------------------------------------------
$ws.Range("A42").Select()
